$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2595.5
$ws.Range("I43").Value = 2495
$ws.Range("J43").Value = 2629
$ws.Range("K43").Value = 2495
$ws.Range("L43").Value = 2629
$ws.Range("M43").Value = -2426
$ws.Range("N43").Value = -2767
$ws.Range("H55").Value = 911.05
$ws.Range("I55").Value = 1111.3125
$ws.Range("J55").Value = 110
$ws.Range("K55").Value = 1111.3125
$ws.Range("L55").Value = 110
$ws.Range("M55").Value = -897.3125
$ws.Range("N55").Value = -538
$ws.Range("H112").Value = 1725784.1
$ws.Range("J112").Value = 2021491.4
$ws.Range("L112").Value = 6064474.199999999
$ws.Range("N112").Value = -6066690.199999999
$ws.Range("H129").Value = 1831.7273
$ws.Range("I129").Value = 805.8333
$ws.Range("K129").Value = 2417.4999
$ws.Range("M129").Value = 2582.5001
$ws.Range("H137").Value = 8452.056
$ws.Range("I137").Value = 11498.091
$ws.Range("K137").Value = 34494.273
$ws.Range("M137").Value = -31944.273
$ws.Range("H138").Value = 1982.8081
$ws.Range("I138").Value = 1607.6316
$ws.Range("J138").Value = 2216.5247
$ws.Range("K138").Value = 4822.8948
$ws.Range("L138").Value = 6649.5741
$ws.Range("M138").Value = 317.1052
$ws.Range("N138").Value = -16929.5741
$ws.Range("H141").Value = 5173.25
$ws.Range("I141").Value = 5173.25
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 15519.75
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -10339.75
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4573.276
$ws.Range("I2").Value = 4768.269
$ws.Range("K2").Value = 4768.269
$ws.Range("M2").Value = -4655.269
$ws.Range("H23").Value = 22999
$ws.Range("J23").Value = 22999
$ws.Range("L23").Value = 22999
$ws.Range("N23").Value = -23517
$ws.Range("H32").Value = 15848.49
$ws.Range("I32").Value = 14422.617
$ws.Range("K32").Value = 14422.617
$ws.Range("M32").Value = -14135.617
$ws.Range("H88").Value = 2214
$ws.Range("I88").Value = 1992.1666
$ws.Range("J88").Value = 2361.889
$ws.Range("K88").Value = 1992.1666
$ws.Range("L88").Value = 2361.889
$ws.Range("M88").Value = -1586.1666
$ws.Range("N88").Value = -3173.889
$ws.Range("H91").Value = 2214
$ws.Range("I91").Value = 1992.1666
$ws.Range("J91").Value = 2361.889
$ws.Range("K91").Value = 1992.1666
$ws.Range("L91").Value = 2361.889
$ws.Range("M91").Value = -588.1666
$ws.Range("N91").Value = -5169.889
$ws.Range("H116").Value = 4573.276
$ws.Range("I116").Value = 4768.269
$ws.Range("K116").Value = 4768.269
$ws.Range("M116").Value = -2474.269
$ws.Range("H125").Value = 20000
$ws.Range("J125").Value = 20000
$ws.Range("L125").Value = 20000
$ws.Range("N125").Value = -29840

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4573.276
$ws.Range("I3").Value = 4768.269
$ws.Range("K3").Value = 4768.269
$ws.Range("M3").Value = -4654.269
$ws.Range("H20").Value = 6602.3335
$ws.Range("I20").Value = 4808
$ws.Range("K20").Value = 4808
$ws.Range("M20").Value = -4561
$ws.Range("H54").Value = 6749.5
$ws.Range("I54").Value = 6749.5
$ws.Range("K54").Value = 6749.5
$ws.Range("M54").Value = -6265.5
$ws.Range("H99").Value = 10057.647
$ws.Range("J99").Value = 5312.5
$ws.Range("L99").Value = 5312.5
$ws.Range("N99").Value = -8308.5
$ws.Range("H107").Value = 4712
$ws.Range("J107").Value = 4712
$ws.Range("L107").Value = 4712
$ws.Range("N107").Value = -8552
$ws.Range("H134").Value = 1584.2839
$ws.Range("I134").Value = 1280.45
$ws.Range("K134").Value = 3841.35
$ws.Range("M134").Value = -1306.35

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1373.5555
$ws.Range("I5").Value = 258
$ws.Range("J5").Value = 2266
$ws.Range("K5").Value = 258
$ws.Range("L5").Value = 2266
$ws.Range("M5").Value = -146
$ws.Range("N5").Value = -2490
$ws.Range("H8").Value = 251700
$ws.Range("J8").Value = 500000
$ws.Range("L8").Value = 500000
$ws.Range("N8").Value = -500280
$ws.Range("H31").Value = 26339420
$ws.Range("I31").Value = 2229.3333
$ws.Range("K31").Value = 2229.3333
$ws.Range("M31").Value = -1934.3333
$ws.Range("H34").Value = 26339420
$ws.Range("I34").Value = 2229.3333
$ws.Range("K34").Value = 2229.3333
$ws.Range("M34").Value = -2027.3333
$ws.Range("H99").Value = 3017.1177
$ws.Range("I99").Value = 2879.3572
$ws.Range("K99").Value = 2879.3572
$ws.Range("M99").Value = -1381.3572
$ws.Range("H126").Value = 3017.1177
$ws.Range("I126").Value = 2879.3572
$ws.Range("K126").Value = 8638.071599999999
$ws.Range("M126").Value = -6168.071599999999
$ws.Range("H132").Value = 3261.75
$ws.Range("I132").Value = 2732.2856
$ws.Range("K132").Value = 8196.856800000001
$ws.Range("M132").Value = -5666.856800000001
$ws.Range("H134").Value = 4142.275
$ws.Range("I134").Value = 4204.9443
$ws.Range("K134").Value = 12614.8329
$ws.Range("M134").Value = -10079.8329

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1000000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H131").Value = 7144463.5
$ws.Range("I131").Value = 55556012
$ws.Range("K131").Value = 166668036
$ws.Range("M131").Value = -166662996
$ws.Range("H134").Value = 16521.182
$ws.Range("I134").Value = 16521.182
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 49563.546
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -44493.546
$ws.Range("N134").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 333335070
$ws.Range("I70").Value = 2600
$ws.Range("J70").Value = 1000000000
$ws.Range("K70").Value = 2600
$ws.Range("L70").Value = 1000000000
$ws.Range("M70").Value = -2330
$ws.Range("N70").Value = -1000000540
$ws.Range("H73").Value = 333335070
$ws.Range("I73").Value = 2600
$ws.Range("J73").Value = 1000000000
$ws.Range("K73").Value = 2600
$ws.Range("L73").Value = 1000000000
$ws.Range("M73").Value = -1664
$ws.Range("N73").Value = -1000001872
$ws.Range("H80").Value = 3320.2222
$ws.Range("I80").Value = 3568.8572
$ws.Range("K80").Value = 3568.8572
$ws.Range("M80").Value = -2570.8572
$ws.Range("H83").Value = 3320.2222
$ws.Range("I83").Value = 3568.8572
$ws.Range("K83").Value = 17844.286
$ws.Range("M83").Value = -12852.286
$ws.Range("H132").Value = 30776.324
$ws.Range("I132").Value = 32387.727
$ws.Range("K132").Value = 97163.181
$ws.Range("M132").Value = -94633.181
$ws.Range("H136").Value = 35263.566
$ws.Range("J136").Value = 35263.566
$ws.Range("L136").Value = 105790.698
$ws.Range("N136").Value = -110890.698

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 4864.1177
$ws.Range("I55").Value = 270.3158
$ws.Range("J55").Value = 10682.934
$ws.Range("K55").Value = 270.3158
$ws.Range("L55").Value = 10682.934
$ws.Range("M55").Value = -97.31580000000002
$ws.Range("N55").Value = -11028.934
$ws.Range("H68").Value = 2053.9666
$ws.Range("I68").Value = 2011.8276
$ws.Range("J68").Value = 3276
$ws.Range("K68").Value = 2011.8276
$ws.Range("L68").Value = 3276
$ws.Range("M68").Value = -1262.8276
$ws.Range("N68").Value = -4774
$ws.Range("H71").Value = 2053.9666
$ws.Range("I71").Value = 2011.8276
$ws.Range("J71").Value = 3276
$ws.Range("K71").Value = 10059.138
$ws.Range("L71").Value = 16380
$ws.Range("M71").Value = -6315.138000000001
$ws.Range("N71").Value = -23868
$ws.Range("H136").Value = 3317.36
$ws.Range("I136").Value = 2518.2354
$ws.Range("K136").Value = 7554.706200000001
$ws.Range("M136").Value = -5004.706200000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2323.9
$ws.Range("I126").Value = 1301.4073
$ws.Range("K126").Value = 3904.2219
$ws.Range("M126").Value = -1434.2219
$ws.Range("H132").Value = 10527906
$ws.Range("I132").Value = 18182666
$ws.Range("J132").Value = 2609.625
$ws.Range("K132").Value = 54547998
$ws.Range("L132").Value = 7828.875
$ws.Range("M132").Value = -54545468
$ws.Range("N132").Value = -12888.875
$ws.Range("H135").Value = 56125
$ws.Range("J135").Value = 58428.57
$ws.Range("L135").Value = 58428.57
$ws.Range("N135").Value = -68568.57000000001
$ws.Range("H136").Value = 3606.4614
$ws.Range("I136").Value = 3297.6177
$ws.Range("J136").Value = 5706.6
$ws.Range("K136").Value = 9892.8531
$ws.Range("L136").Value = 17119.8
$ws.Range("M136").Value = -7342.8531
$ws.Range("N136").Value = -22219.8
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("H141").Value = 10000
$ws.Range("I141").Value = 10000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 10000
$ws.Range("L141").Value = 10000
$ws.Range("M141").Value = -4820
$ws.Range("N141").ClearContents()
